$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing the existing separator/total/paid/credit rows
# (old 16-21) down to 17-22.
$ws.Rows.Item(16).Insert()

# Re-use the exact formatting of the row above (dates/text/hours columns) for
# the freshly inserted row instead of the fresh blank-row styles Insert()
# hands us, so the new cells land on the same style records as row 15.
$ws.Range("B15:F15").Copy()
$ws.Range("B16:F16").PasteSpecial(-4122) # xlPasteFormats

# Populate the new data row: dates, work description and hours, matching the
# new shared-string entry "Data extracted from new prodrome pdfs".
$ws.Cells.Item(16, 2).Value = 42919
$ws.Cells.Item(16, 3).Value = 42925
$ws.Cells.Item(16, 4).Value = "Data extracted from new prodrome pdfs"
$ws.Cells.Item(16, 5).Value = 5

# The TOTAL HOURS sum (now on row 18) needs to include the new row 16.
$ws.Range("E18").Formula = "=SUM(E5:E16)"

# Conditional formatting ranges don't auto-shift with the row insert in this
# host, so retarget them to the rows they now land on (19->20, 21->22).
$ws.Range("E19").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E20"))
$ws.Range("E21").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E22"))
